$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '258.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.01%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.06%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.217'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.44%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05925'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.72%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.706'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.88%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8666'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.31%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.006'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '14.00%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1417'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.09%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07178'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.37%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03152'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.17%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09219'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.17%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001538'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.20%'

$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006075'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.45%'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005915'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.22%'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.498'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.09%'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.269'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.66%'

$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.207'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.58%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3143'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.79%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '6.33%'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.27%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.539'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.44%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04185'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.05%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.44%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001216'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.58%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004527'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '8.57%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.02%'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001483'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '2.67%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03830'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.95%'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '15.66%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1106'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.64%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002359'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '7.29%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01083'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '14.41%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005426'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.70%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.01%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '22.38%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002232'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.02%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.01%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
